# Applies the cryptos.xlsx price/volume/coin-row updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each touched cell to remain plain text (matches the original inlineStr/text
# cells) so Excel does not silently reinterpret numeric-looking strings (e.g. "0.0950")
# as numbers and drop significant trailing zeros.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '34.398.49'
Set-TextValue 'E2' '  +0.94%  '
Set-TextValue 'D3' '1.793.80'
Set-TextValue 'E3' '  +0.54%  '
Set-TextValue 'E4' '  -0.34%  '
Set-TextValue 'D5' '226.56'
Set-TextValue 'E5' '  +0.44%  '
Set-TextValue 'E6' '  +1.67%  '
Set-TextValue 'E7' '  -0.23%  '
Set-TextValue 'D8' '32.44'
Set-TextValue 'E8' '  +1.99%  '
Set-TextValue 'E9' '  +1.72%  '
Set-TextValue 'D10' '0.0692'
Set-TextValue 'E10' '  +0.72%  '
Set-TextValue 'D11' '0.0950'
Set-TextValue 'E11' '  +0.76%  '
Set-TextValue 'D12' '2.051.21'
Set-TextValue 'E12' '  +0.36%  '
Set-TextValue 'B13' 'WrappedEther'
Set-TextValue 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D13' '1.793.64'
Set-TextValue 'E13' '  +0.47%  '
Set-TextValue 'B14' 'Chainlink'
Set-TextValue 'C14' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D14' '11.04'
Set-TextValue 'E14' '  -1.06%  '
Set-TextValue 'D15' '0.630'
Set-TextValue 'E15' '  +2.19%  '
Set-TextValue 'D16' '34.347.16'
Set-TextValue 'E16' '  +0.78%  '
Set-TextValue 'D17' '4.22'
Set-TextValue 'E17' '  +0.92%  '
Set-TextValue 'D18' '68.13'
Set-TextValue 'E18' '  +0.44%  '
Set-TextValue 'D19' '0.0₃0801'
Set-TextValue 'E19' '  +3.28%  '
Set-TextValue 'D20' '246.68'
Set-TextValue 'E20' '  +0.59%  '
Set-TextValue 'D21' '10.96'
Set-TextValue 'E21' '  +1.77%  '
Set-TextValue 'E22' '  -0.13%  '
Set-TextValue 'E23' '  +2.12%  '
Set-TextValue 'E24' '  +0.93%  '
Set-TextValue 'D25' '162.39'
Set-TextValue 'E25' '  +0.67%  '
Set-TextValue 'D26' '7.19'
Set-TextValue 'E26' '  +1.17%  '
Set-TextValue 'D27' '16.39'
Set-TextValue 'E27' '  +0.78%  '
Set-TextValue 'E28' '  +2.22%  '
Set-TextValue 'E29' '  -0.22%  '
Set-TextValue 'D30' '3.94'
Set-TextValue 'E30' '  +10.07%  '
Set-TextValue 'E31' '  +0.18%  '
Set-TextValue 'E32' '  +0.80%  '
Set-TextValue 'E33' '  +3.88%  '
Set-TextValue 'D34' '1.83'
Set-TextValue 'E34' '  +1.31%  '
Set-TextValue 'D35' '1.441.60'
Set-TextValue 'E35' '  -0.96%  '
Set-TextValue 'D36' '2.62'
Set-TextValue 'E36' '  +9.93%  '
Set-TextValue 'D37' '0.663'
Set-TextValue 'E37' '  +3.38%  '
Set-TextValue 'D38' '1.05'
Set-TextValue 'E38' '  +1.92%  '
Set-TextValue 'E39' '  -1.14%  '
Set-TextValue 'D40' '83.18'
Set-TextValue 'E40' '  +4.40%  '
Set-TextValue 'D41' '2.39'
Set-TextValue 'E41' '  +0.83%  '
Set-TextValue 'D42' '14.01'
Set-TextValue 'E42' '  +5.42%  '
Set-TextValue 'B43' 'MXToken'
Set-TextValue 'C43' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D43' '2.75'
Set-TextValue 'E43' '  +0.43%  '
Set-TextValue 'B44' 'ARBITRUM'
Set-TextValue 'C44' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D44' '0.930'
Set-TextValue 'E44' '  +1.76%  '
Set-TextValue 'D45' '0.0520'
Set-TextValue 'E45' '  +2.16%  '
Set-TextValue 'D46' '6.07'
Set-TextValue 'E46' '  +0.39%  '
Set-TextValue 'E47' '  -0.19%  '
Set-TextValue 'D48' '1.946.45'
Set-TextValue 'E48' '  +0.08%  '
Set-TextValue 'D49' '105.49'
Set-TextValue 'E49' '  -1.49%  '
Set-TextValue 'B50' 'PaxDollar'
Set-TextValue 'C50' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D50' '1.00'
Set-TextValue 'E50' '  -0.12%  '
Set-TextValue 'B51' 'BabyDogeCoin'
Set-TextValue 'C51' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D51' '0.0₆0130'
Set-TextValue 'E51' '  -5.51%  '
